$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 23 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A23").Value = "Wanneer zijn jullie open?"
$logs.Range("B23").Value = "mailmind.test@zohomail.eu"
$logs.Range("C23").Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D23").Value = "Openingstijden / Locatie"
$logs.Range("E23").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F23").Value = "2025-06-26 22:15:54"
$logs.Range("G23").Value = "Ja"
$logs.Range("H23").Value = "Nee"
$logs.Range("I23").Value = "Ja"

# Multi-line text in E23 triggers an implicit wrap-based row autofit that
# would otherwise pin an explicit ht/customHeight on row 23; re-autofitting
# the row clears that back to the sheet's default (matches source rows,
# none of which carry an explicit row height).
$logs.Rows.Item(23).EntireRow.AutoFit() | Out-Null

# --- Extend conditional formatting ranges to include the new row 23 ---
$cfCols = @("D", "G", "H", "I")
foreach ($col in $cfCols) {
    $startCell = $col + "2"
    $rangeAddr = $col + "2:" + $col + "23"
    $fcs = $logs.Range($startCell).FormatConditions
    $newRange = $logs.Range($rangeAddr)
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append row 6 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Openingstijden / Locatie"
$dash.Range("B6").Value = 1

# --- Update chart series ranges to include the new Dashboard row 6 ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$ser.Values = "='Dashboard'!`$B`$2:`$B`$6"
